$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stash the current header formatting (bold/fill/border/alignment) in a
# scratch cell outside the used range so it can be restored after the
# table is created (creating a ListObject/table on pre-formatted header
# cells would otherwise bake a header-row dxf into styles.xml).
$scratch = $ws.Range("W1")
$ws.Range("A1").Copy()
$scratch.PasteSpecial(-4122)   # xlPasteFormats

$headerRange = $ws.Range("A1:U1")
$headerRange.ClearFormats()

# Rename header cells: "_old" -> "_FV2304" and "_new" -> "_FV2310"
$fv2304 = @("Segmentname_FV2304","Segmentgruppe_FV2304","Segment_FV2304","Datenelement_FV2304","Segment ID_FV2304","Code_FV2304","Qualifier_FV2304","Beschreibung_FV2304","Bedingungsausdruck_FV2304","Bedingung_FV2304")
$fv2310 = @("Segmentname_FV2310","Segmentgruppe_FV2310","Segment_FV2310","Datenelement_FV2310","Segment ID_FV2310","Code_FV2310","Qualifier_FV2310","Beschreibung_FV2310","Bedingungsausdruck_FV2310","Bedingung_FV2310")

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2304[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2310[$i]
}

# Turn the used range into an Excel table (ListObject)
$rng = $ws.Range("A1:U61")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Restore the original header formatting now that the table exists.
$scratch.Copy()
$headerRange.PasteSpecial(-4122)   # xlPasteFormats
$scratch.Clear()

# Freeze the header row (split after row 1)
$ws.Range("A2").Select()
[void]($excel.ActiveWindow.FreezePanes = $true)
